# Add a Japanese translation column (G) for the English HPO terms in column B.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$translations = @{
    2  = "乳輪の充実"
    3  = "低形成乳輪"
    4  = "乳房肥大"
    5  = "乳房無形成"
    6  = "両側乳房低形成"
    7  = "片側乳房低形成"
    8  = "反転乳頭"
    9  = "過剰乳頭"
    10 = "乾燥乳頭"
    11 = "広い乳頭間距離"
    12 = "短い乳頭間距離"
    13 = "低位乳頭"
    14 = "欠如乳頭"
    15 = "低形成乳頭"
    16 = "隆起乳頭"
    17 = "乳房の非対称性"
    18 = "多巣性乳癌"
    19 = "非浸潤性乳管がん"
    20 = "非浸潤性小葉がん"
    21 = "乳房線維腺腫"
    22 = "乳房腫瘤"
    23 = "女性化乳房"
    24 = "乳汁漏出症"
    25 = "無月経"
}

foreach ($row in 2..25) {
    $ws.Cells.Item($row, 7).Value = $translations[$row]
}

$ws.Rows("1:1").Select()
